# "performance improvement + demi palettes"
# The METRO pickup point is discontinued (its row is removed entirely,
# shifting every following row up by one), and Carrefour Supply Chain's
# weight per pickup ("Poids par ramasse(kg)") increases from 250 to 800
# (half-pallets instead of quarter-pallets).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the METRO row (row 2); Leclerc/Auchan/Super U/Carrefour rows
# all shift up one position.
$ws.Rows(2).Delete()

# Carrefour Supply Chain is now the last data row (row 7) after the
# shift; bump its pickup weight from 250kg to 800kg.
$ws.Range("G7").Value = 800

# Carrefour's "Poids" cell previously had mismatched formatting (no
# border/general alignment, same as the now-removed METRO row) while
# every other pickup point in the column is bordered & right-aligned.
# Normalize it to match the rest of the column now that it is the
# last row.
$ws.Range("G2").Copy()
$ws.Range("G7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
